# LSTM_mv.xlsx data refresh — appends the latest rows of market data
# (dates 2023-12-14 .. 2023-12-22, i.e. serials 45280-45295) to the
# per-pair worksheets, replaces the "Nan" placeholders on D5_EUR with
# the now-known actuals, and updates the trailing selection on each
# touched sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# D1_USD  (sheet1) — rows 88-96 new, row 97 gets its lone C value
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("D1_USD")
$ws.Activate()

$rows = @(
  @(88, 45280, 3.9370440000000002, 4.0000270000000002),
  @(89, 45281, 3.9678070000000001, 3.9410834000000001),
  @(90, 45282, 3.9279000000000002, 3.9743461999999998),
  @(91, 45286, 3.9323109999999999, 3.8880050000000002),
  @(92, 45287, 3.9164859999999999, 3.9465865999999998),
  @(93, 45288, 3.886606,           3.9129738999999999),
  @(94, 45289, 3.9155280000000001, 3.8855276000000001),
  @(95, 45293, 3.9369000000000001, 3.9294790000000002),
  @(96, 45294, 3.9897559999999999, 3.9566110000000001)
)

foreach ($row in $rows) {
  $r = $row[0]
  $ws.Range("A$r").Value = $row[1]
  $ws.Range("B$r").Value = $row[2]
  $ws.Range("C$r").Value = $row[3]
  $ws.Range("D$r").Formula = "=B$r-C$r"
  $ws.Range("E$r").Formula = "=IF(D$r<0,1,0)"
}
$ws.Range("C97").Value = 4.0088629999999998

# copy number formats down from the last fully-formatted row (87)
$src = $ws.Range("A87:E87")
$dst = $ws.Range("A88:E96")
$src.Copy()
$dst.PasteSpecial(-4122)

$ws.Range("A97").Select()

# ---------------------------------------------------------------
# D1_EUR  (sheet3) — rows 362-370 new, row 371 gets its lone C value
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("D1_EUR")
$ws.Activate()

$rows = @(
  @(362, 45280, 4.32315,             4.3047713999999999),
  @(363, 45281, 4.34213,             4.300592),
  @(364, 45282, 4.3228,              4.3351129999999998),
  @(365, 45286, 4.3332100000000002,  4.3285612999999996),
  @(366, 45287, 4.3247400000000003,  4.3349209999999996),
  @(367, 45288, 4.3173000000000004,  4.3282090000000002),
  @(368, 45289, 4.3335499999999998,  4.3232400000000002),
  @(369, 45293, 4.3479000000000001,  4.3370059999999997),
  @(370, 45294, 4.3661899999999996,  4.3579197000000001)
)

foreach ($row in $rows) {
  $r = $row[0]
  $ws.Range("A$r").Value = $row[1]
  $ws.Range("B$r").Value = $row[2]
  $ws.Range("C$r").Value = $row[3]
  $ws.Range("D$r").Formula = "=B$r-C$r"
  $ws.Range("E$r").Formula = "=IF(D$r<0,1,0)"
}
$ws.Range("C371").Value = 4.4089565000000004

$src = $ws.Range("A361:E361")
$dst = $ws.Range("A362:E370")
$src.Copy()
$dst.PasteSpecial(-4122)

$ws.Range("A370:B370").Select()

# ---------------------------------------------------------------
# D5_EUR  (sheet5) — B55:B59 "Nan" placeholders now have real
# predictions, plus brand-new rows 60-64 (row 64 keeps a "Nan" in B)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("D5_EUR")
$ws.Activate()

$ws.Range("B55").Value = 4.32315
$ws.Range("B56").Value = 4.34213
$ws.Range("B57").Value = 4.3228
$ws.Range("B58").Value = 4.3332100000000002
$ws.Range("B59").Value = 4.3247400000000003

$rows = @(
  @(60, 45288, 4.3173000000000004, 4.2615129999999999),
  @(61, 45289, 4.3335499999999998, 4.2706203),
  @(62, 45293, 4.33988,            4.2816386),
  @(63, 45294, 4.3661899999999996, 4.2987776000000002)
)
foreach ($row in $rows) {
  $r = $row[0]
  $ws.Range("A$r").Value = $row[1]
  $ws.Range("B$r").Value = $row[2]
  $ws.Range("C$r").Value = $row[3]
}
$ws.Range("A64").Value = 45295
$ws.Range("B64").Value = "Nan"
$ws.Range("C64").Value = 4.3172812

$src = $ws.Range("A59:C59")
$dst = $ws.Range("A60:C64")
$src.Copy()
$dst.PasteSpecial(-4122)

$ws.Range("B65").Select()

# ---------------------------------------------------------------
# D1_OIL  (sheet6) — rows 34-41 new (last two rows only have A/B)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("D1_OIL")
$ws.Activate()

$fullRows = @(
  @(34, 45280, 73.440002000000007, 71.343400000000003),
  @(35, 45281, 73.889999000000003, 70.639300000000006),
  @(36, 45282, 73.559997999999993, 71.102999999999994),
  @(37, 45286, 73.559997999999993, 67.181200000000004),
  @(38, 45287, 74.110000999999997, 67.802999999999997),
  @(39, 45288, 71.769997000000004, 68.910200000000003)
)
foreach ($row in $fullRows) {
  $r = $row[0]
  $ws.Range("A$r").Value = $row[1]
  $ws.Range("B$r").Value = $row[2]
  $ws.Range("C$r").Value = $row[3]
  $ws.Range("D$r").Formula = "=B$r-C$r"
  $ws.Range("E$r").Formula = "=D$r/C$r"
}

$src = $ws.Range("A33:E33")
$dst = $ws.Range("A34:E39")
$src.Copy()
$dst.PasteSpecial(-4122)

$ws.Range("A40").Value = 45289
$ws.Range("B40").Value = 71.650002000000001
$ws.Range("A41").Value = 45293
$ws.Range("B41").Value = 70.379997000000003

$src2 = $ws.Range("A39:B39")
$dst2 = $ws.Range("A40:B41")
$src2.Copy()
$dst2.PasteSpecial(-4122)

$ws.Range("C40").Select()

# ---------------------------------------------------------------
# Leave D5_EUR as the active sheet/tab, matching the saved workbook
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("D5_EUR")
$ws.Activate()
$ws.Range("B65").Select()
